$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted at row 189; all the existing
# rows from 189 down to 202 shift down by one (to 190..203).
$ws.Rows("189").Insert()

$ws.Cells.Item(189, 1).Value = 7
$ws.Cells.Item(189, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(189, 3).Value = "Ñuble"
$ws.Cells.Item(189, 4).Value = 44585
$ws.Cells.Item(189, 5).Value = 16
$ws.Cells.Item(189, 6).Value = 100112009
$ws.Cells.Item(189, 7).Value = "Acelga"
$ws.Cells.Item(189, 8).Value = "Sin especificar"
$ws.Cells.Item(189, 9).Value = "Primera"
$ws.Cells.Item(189, 10).Value = 60
$ws.Cells.Item(189, 11).Value = 350
$ws.Cells.Item(189, 12).Value = 400
$ws.Cells.Item(189, 13).Value = 375
$ws.Cells.Item(189, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(189, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(189, 16).Value = 375
$ws.Cells.Item(189, 17).Value = 1
$ws.Cells.Item(189, 18).Value = "Hortaliza"
